$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the existing "text" cell style (same as A2/B2 already had) to the
# other string-valued cells in row 2 before writing their values, so the
# values land as shared strings under the same style index.
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# Populate row 2 with server data (order chosen so new shared strings are
# appended in the same order as the target workbook: IP, then the server
# name, then the server id).
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "LoginServer_1"
$ws.Range("C2").Value = "LoginServer_1"
$ws.Range("B2").Value = "000106001"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 6001

# Restrict the list data validation to start from row 3, since row 2 now
# holds explicit server data instead of being part of the dropdown range.
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update the active selection to reflect where the user last clicked.
$ws.Range("G5").Select()
